# edit.ps1 - apply the "adjustments to functional and nonfunctional req" commit
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "1.2 Throughput:" note -> split out a bold/blue " NOTE " run
#    " (Check the average of other similar website; applications)"
#    becomes " (" + " NOTE " + "Check the average of other similar website; applications)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " (Check the average of other similar website; applications)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " ( NOTE Check the average of other similar website; applications)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "1.4 Reliability:" note -> split out a "NOTE " run
#    "(Recovery time also based on other applications)"
#    becomes "(" + "NOTE " + "Recovery time also based on other applications)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(Recovery time also based on other applications)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(NOTE Recovery time also based on other applications)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "2.1 User interface:" paragraph rewrite
#    " (separate interface for each different type of user)" + ", and accessible, "
#    becomes
#    " by separating the interfaces for different types of users, and accessible."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " (separate interface for each different type of user), and accessible, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " by separating the interfaces for different types of users, and accessible.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Mobile application note -> split out a " NOTE " run
#    "(Mobile Application to inc usability and be more usable for the user "
#    becomes "(" + " NOTE " + "Mobile Application to inc usability and be more usable for the user "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(Mobile Application to inc usability and be more usable for the user ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "( NOTE Mobile Application to inc usability and be more usable for the user ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Drop the stray lastRenderedPageBreak before "Requirement 3: Security"
#    (re-stamping the run's text drops the stale pagination hint)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Requirement 3: Security",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Requirement 3: Security", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "3.1 Authentication:" note -> split out a "  NOTE " run and trim the
#    trailing clause
#    " (Optional phone number or Email and add this info in the functional requirements)."
#    becomes " (" + "  NOTE " + "Optional phone number or Email)."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " (Optional phone number or Email and add this info in the functional requirements).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (  NOTE Optional phone number or Email).", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) "3.2 Authorization:" paragraph - "has to be" -> "must be"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Access to different services and features has to be controlled based on user roles to ensure that each user has access to the allowed services only.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Access to different services and features must be controlled based on user roles to ensure that each user has access to the allowed services only.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Insert a new blank paragraph after "Maintain integrity of user data..."
#    (inherits the ind=1440 / sz=26 formatting of that paragraph)
# ---------------------------------------------------------------------------
$rngIntegrity = $d.Content
$rngIntegrity.Find.Execute("Maintain integrity of user data by ensuring the accuracy, consistency, and reliability of data in the system.") | Out-Null
$rngIntegrity.Collapse(0) | Out-Null
$rngIntegrity.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 9) "4.1 Modularity:" - drop the "(we will adopt...)" aside entirely
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Modularity should be taken into consideration when building the system architecture, for ease in updates and additions to the code (we will adopt the service oriented architecture to achieve modularity).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Modularity should be taken into consideration when building the system architecture, for ease in updates and additions to the code.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 10-14) Restructure the tail of the document:
#   - "4.2 Code maintainability: (Delete)" -> "4.2 Documentation: "
#   - remove the old "Coding standards..." paragraph
#   - remove the old "4.3 Documentation: " heading paragraph
#   - keep "Comprehensive documentation..." paragraph as-is
#   - keep the blank paragraph after it
#   - remove everything from "Requirement 5: Monitoring" onward
# ---------------------------------------------------------------------------

# 10) Re-text the "4.2 Code maintainability: (Delete)" heading paragraph
$d.Content.Find.Execute(
    "4.2 Code maintainability: (Delete)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "4.2 Documentation: ", 2) | Out-Null

# 11+12) Delete the "Coding standards..." paragraph and the old
#        "4.3 Documentation: " heading paragraph (both now sit between the
#        rewritten 4.2 heading and "Comprehensive documentation...")
$rngHeading42 = $d.Content
$rngHeading42.Find.Execute("4.2 Documentation:") | Out-Null
$paraHeading42 = $rngHeading42.Paragraphs(1)

$rngHeading43 = $d.Content
$rngHeading43.Find.Execute("4.3 Documentation:") | Out-Null
$paraHeading43 = $rngHeading43.Paragraphs(1)

$midDelStart = $paraHeading42.Range.End
$midDelEnd = $paraHeading43.Range.End
$d.Range($midDelStart, $midDelEnd).Delete() | Out-Null

# 14) Remove "Requirement 5: Monitoring" and everything after it, but keep
#     the blank paragraph that follows "Comprehensive documentation..."
$rngComprehensive = $d.Content
$rngComprehensive.Find.Execute("Comprehensive documentation is needed to ensure ease of future development.") | Out-Null
$paraComprehensive = $rngComprehensive.Paragraphs(1)
$paraBlank = $paraComprehensive.Next()
$tailDelStart = $paraBlank.Range.End
$tailDelEnd = $d.Content.End
$d.Range($tailDelStart, $tailDelEnd).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 15) Add triple page borders to the (only) section
# ---------------------------------------------------------------------------
$sectionBorders = $d.Sections(1).Borders
$sectionBorders.Item(1).LineStyle = 8
$sectionBorders.DistanceFromTop = 24
$sectionBorders.DistanceFromBottom = 24
$sectionBorders.DistanceFromLeft = 24
$sectionBorders.DistanceFromRight = 24

Write-Host "done part 1"
